# style: Fix spelling mistake in excel template #71
#
# The shared string used by cell A28 of the "Hoja1" worksheet was
# missing the final "n" of "opinión": "Se utilizan encuestas de
# opinió" -> "Se utilizan encuestas de opinión".
#
# Re-assigning the cell's Value updates the existing shared-string
# table entry in place (same <si> slot), so no other cell is affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A28").Value = "Se utilizan encuestas de opinión"

# Leave the selection where the user ended up after correcting the
# text and moving on to the next row.
[void]$ws.Range("A29").Select()
